$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'GyXSr256'
$ws.Range("B2").Value = 23110945
$ws.Range("C2").Value = 'auuzkwk41'
$ws.Range("D2").Value = 'Z&sR6$v2'
$ws.Range("F2").Value = 'tSqTiMSy'
$ws.Range("G2").Value = 'zCeb'

$ws.Rows.Item(3).Delete()

$ws.Range("A1:H2").Select() | Out-Null
